$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original style of the Price column, then force text format
# so that numeric-looking strings (e.g. "75.40") are stored as text
# instead of being auto-converted into numbers by Excel.
$origStyleD = $ws.Range("D2:D51").Style
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.121.09"
$ws.Range("E2").Value = "  +4.98%  "
$ws.Range("D3").Value = "2.235.46"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "245.96"
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("D7").Value = "75.40"
$ws.Range("E7").Value = "  +8.35%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  +6.87%  "
$ws.Range("D10").Value = "41.18"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "55.73"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "6.95"
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "2.570.68"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").Value = "14.75"
$ws.Range("E16").Value = "  +6.69%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "0.814"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.224.33"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").Value = "43.004.41"
$ws.Range("E19").Value = "  +5.16%  "
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").Value = "70.96"
$ws.Range("D22").Value = "5.97"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  +7.75%  "
$ws.Range("D24").Value = "230.49"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("E25").Value = "  +11.63%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "174.83"
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "37.58"
$ws.Range("E32").Value = "  +22.44%  "
$ws.Range("D33").Value = "20.33"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").Value = "0.0794"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  +4.69%  "
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").Value = "0.110"
$ws.Range("E37").Value = "  +6.80%  "
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "0.0333"
$ws.Range("E39").Value = "  +16.89%  "
$ws.Range("D40").Value = "13.13"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("D43").Value = "0.199"
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("D44").Value = "60.09"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "105.39"
$ws.Range("E45").Value = "  +7.27%  "
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").Value = "0.0991"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").Value = "0.444"
$ws.Range("E48").Value = "  +21.34%  "
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +3.67%  "
$ws.Range("E51").Value = "  +2.48%  "

# Restore the original style on the Price column so formatting/style
# stays identical to before (only the underlying text content changes).
$ws.Range("D2:D51").Style = $origStyleD

